$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update effected-rows count for existing row 104 (400 -> 490) ---
$ws.Cells.Item(104, 9).Value = "490"

# --- Duplicate row 107 (style donor) into new rows 108, 109, 110 ---
$ws.Rows("107:107").Copy()
$ws.Rows("108:108").Insert(-4121)
$ws.Rows("107:107").Copy()
$ws.Rows("109:109").Insert(-4121)
$ws.Rows("107:107").Copy()
$ws.Rows("110:110").Insert(-4121)
$excel.CutCopyMode = $false

# --- Row 108: updel_107 ---
$ws.Cells.Item(108, 1).Value = "updel_107"
$ws.Cells.Item(108, 2).Value = "y"
$ws.Cells.Item(108, 3).Value = "带有向量索引的表更新数值标量字段"
$ws.Cells.Item(108, 4).Value = "Index"
$ws.Cells.Item(108, 5).Value = "vector_index"
$ws.Cells.Item(108, 6).Value = "vector062"
$ws.Cells.Item(108, 7).Value = "vector062_value1"
$ws.Cells.Item(108, 8).Value = "update `$vector062 set amount=1234.1234 where id<10"
$ws.Cells.Item(108, 9).Value = "9"
$ws.Cells.Item(108, 10).Value = "select id,amount from `$vector062 where id<10"
$ws.Cells.Item(108, 11).Value = "src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_107.csv"
$ws.Cells.Item(108, 12).Value = "csv_containsAll"

# --- Row 109: updel_108 ---
$ws.Cells.Item(109, 1).Value = "updel_108"
$ws.Cells.Item(109, 2).Value = "y"
$ws.Cells.Item(109, 3).Value = "带有向量索引的表更新字符标量字段"
$ws.Cells.Item(109, 4).Value = "Index"
$ws.Cells.Item(109, 5).Value = "vector_index"
$ws.Cells.Item(109, 6).Value = "vector062"
$ws.Cells.Item(109, 7).Value = "vector062_value1"
$ws.Cells.Item(109, 8).Value = "update `$vector062 set address='beijing' where id in (10,20,30)"
$ws.Cells.Item(109, 9).Value = "3"
$ws.Cells.Item(109, 10).Value = "select id,address from `$vector062 where id in (10,20,30)"
$ws.Cells.Item(109, 11).Value = "src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_108.csv"
$ws.Cells.Item(109, 12).Value = "csv_containsAll"

# --- Row 110: updel_109 ---
$ws.Cells.Item(110, 1).Value = "updel_109"
$ws.Cells.Item(110, 2).Value = "y"
$ws.Cells.Item(110, 3).Value = "带有向量索引的表更新日期标量字段"
$ws.Cells.Item(110, 4).Value = "Index"
$ws.Cells.Item(110, 5).Value = "vector_index"
$ws.Cells.Item(110, 6).Value = "vector062"
$ws.Cells.Item(110, 7).Value = "vector062_value1"
$ws.Cells.Item(110, 8).Value = "update `$vector062 set birthday='2023-09-20' where id=100"
$ws.Cells.Item(110, 9).Value = "1"
$ws.Cells.Item(110, 10).Value = "select id,birthday from `$vector062 where id=100"
$ws.Cells.Item(110, 11).Value = "src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_109.csv"
$ws.Cells.Item(110, 12).Value = "csv_containsAll"

# --- Fix style for C108/H108 which should use the "s=1" style instead of the copied "s=7" ---
$ws.Cells.Item(104, 3).Copy()
$ws.Cells.Item(108, 3).PasteSpecial(-4122)
$ws.Cells.Item(104, 8).Copy()
$ws.Cells.Item(108, 8).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(108, 3).Value = "带有向量索引的表更新数值标量字段"
$ws.Cells.Item(108, 8).Value = "update `$vector062 set amount=1234.1234 where id<10"

# --- Update the view: selection + top-left cell ---
$ws.Range("G104").Select()
$excel.ActiveWindow.ScrollRow = 82
$excel.ActiveWindow.ScrollColumn = 1

Write-Output "edit complete"